# Adds the "I0" and "IF" columns (I and J) to the save-data sheet, as per
# commit "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): copy formatting (border/bold/alignment) from the
# existing "IP" header cell (H1) onto the two new header cells. ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values for rows 2-73, columns I (I0) and J (IF) ---
$I0 = @(8,6,7,8,7,11,8,5,8,6,8,7,5,9,9,3,7,4,6,9,8,7,8,8,7,9,4,8,9,5,5,6,8,9,8,9,8,7,6,8,6,9,6,8,6,6,9,7,6,9,8,9,6,7,9,9,7,7,8,8,5,6,5,7,6,5,4,5,5,4,4,3)
$IF = @(8,6,7,8,7,11,8,5,8,6,8,7,6,9,9,4,8,5,6,9,8,7,8,8,7,9,4,8,9,5,6,6,8,9,8,9,8,7,6,8,7,9,7,8,6,6,9,8,6,9,8,9,6,8,9,9,7,7,8,8,5,6,7,8,6,5,4,5,5,4,4,3)

$firstRow = 2
$lastRow = 73
$n = $lastRow - $firstRow + 1

$dataI = New-Object 'object[,]' $n,1
$dataJ = New-Object 'object[,]' $n,1
for ($k = 0; $k -lt $n; $k++) {
    $dataI[$k,0] = $I0[$k]
    $dataJ[$k,0] = $IF[$k]
}

$ws.Range("I${firstRow}:I${lastRow}").Value = $dataI
$ws.Range("J${firstRow}:J${lastRow}").Value = $dataJ

Write-Output "I0/IF columns written"
